$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")

# Sprint Backlog 4 (SB4) user stories for frontend development,
# appended to the Product Backlog sheet starting at row 94.
$sb4Data = @(
    ,@(4, 'SB4/US1', 'Developer', 'design the frontend architecture', 'the UI is well-structured', 4, 'HC', 45992, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US2', 'Developer', 'create the HTML structure', 'the UI has proper semantic markup', 4, 'HC', 45995, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US3', 'Developer', 'implement CSS styling', 'the interface is visually appealing', 4, 'HC', 45998, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US4', 'Developer', 'implement JavaScript functionality', 'the frontend interacts with the API', 4, 'HC', 46001, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US5', 'User', 'submit questions', 'I can get answers from the Virtual TA', 4, 'HC', 46004, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US6', 'User', 'upload images', 'the TA can analyze screenshots', 4, 'HC', 46007, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US7', 'User', 'view answers', 'I can get information from the Virtual TA', 4, 'HC', 46010, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US8', 'Developer', 'implement error handling', 'users get helpful feedback', 4, 'HC', 46013, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US9', 'Developer', 'add CORS support', 'the API works with the frontend', 4, 'HC', 46016, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US10', 'Developer', 'create a chat history feature', 'users can review past queries', 4, 'HC', 46019, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US11', 'User', 'the interface to be responsive', 'I can use it on mobile devices', 4, 'HC', 46022, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US12', 'Developer', 'optimize performance', 'the app loads quickly', 4, 'HC', 46025, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US13', 'Developer', 'add accessibility features', 'the app is usable by everyone', 4, 'HC', 46028, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US14', 'Developer', 'create frontend documentation', 'others can understand the code', 4, 'HC', 46031, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US15', 'Developer', 'test the frontend', 'bugs are caught early', 4, 'HC', 46034, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US16', 'Developer', 'deploy the frontend', 'users can access it', 4, 'HC', 46037, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US17', 'Developer', 'add theme support', 'users can customize appearance', 4, 'HC', 46040, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US18', 'Developer', 'integrate API endpoint', 'frontend can communicate with backend', 4, 'HC', 46043, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US19', 'User', 'see example queries', 'I know what to ask', 4, 'HC', 46046, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US20', 'Developer', 'add analytics tracking', 'usage can be monitored', 4, 'HC', 46049, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US21', 'Developer', 'create a landing page', 'users understand the Virtual TA', 4, 'HC', 46052, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US22', 'Developer', 'add loading animations', 'users know the app is working', 4, 'HC', 46055, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US23', 'User', 'copy answers', 'I can save useful information', 4, 'HC', 46058, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US24', 'Developer', 'implement rate limiting', 'API is not overwhelmed', 4, 'HC', 46061, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US25', 'Developer', 'add final polish', 'the frontend is production-ready', 4, 'HC', 46064, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US26', 'Developer', 'update project documentation', 'Sprint 4 work is captured', 4, 'HC', 46067, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US27', 'Developer', 'create user guide', 'users know how to use the Virtual TA', 4, 'HC', 46070, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US28', 'Developer', 'implement feedback mechanism', 'users can report issues', 4, 'HC', 46073, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US29', 'Developer', 'configure CI/CD', 'frontend deploys automatically', 4, 'HC', 46076, 'Harsh Choudhary', 'HC')
    ,@(4, 'SB4/US30', 'Developer', 'conduct final review', 'Sprint 4 is complete', 4, 'HC', 46079, 'Harsh Choudhary', 'HC')
)

$startRow = 94
for ($i = 0; $i -lt $sb4Data.Count; $i++) {
    $r = $startRow + $i
    $item = $sb4Data[$i]
    $ws.Cells.Item($r, 1).Value = $item[0]   # A: Sprint Backlog #
    $ws.Cells.Item($r, 2).Value = $item[1]   # B: US ID
    $ws.Cells.Item($r, 3).Value = $item[2]   # C: As a/an
    $ws.Cells.Item($r, 4).Value = $item[3]   # D: I want to
    $ws.Cells.Item($r, 5).Value = $item[4]   # E: So that
    $ws.Cells.Item($r, 6).Value = $item[5]   # F: Priority
    $ws.Cells.Item($r, 7).Value = $item[6]   # G: Responsible
    $ws.Cells.Item($r, 8).Value = $item[7]   # H: Estimate date (serial)
    $ws.Cells.Item($r, 8).NumberFormat = "yyyy-mm-dd h:mm:ss"
    $ws.Cells.Item($r, 11).Value = $item[8]  # K: Student Name
    $ws.Cells.Item($r, 12).Value = $item[9]  # L: Code
}

